$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that sat alone in the empty paragraph right
#    after "Reviewer #2: MAJOR COMMENTS" (paragraph becomes truly empty).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Replace the "xxxxxxxxxxxxxxxxxxxxxxxx" placeholder answer that directly
#    follows the `"Submitted" papers cannot be cited.` remark with the real
#    response text "They were removed." Using InsertXML on the paragraph's
#    own Range keeps the existing paragraph/run formatting (blue 0066CC)
#    while avoiding leftover <w:proofErr/> spell-check markers that a plain
#    Range.Text edit would leave behind.
# ---------------------------------------------------------------------------
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("`"Submitted`" papers cannot be cited.") | Out-Null
$afterAnchor = $anchor.End

$searchScope = $d.Range($afterAnchor, $afterAnchor + 400)
$searchScope.Find.Execute("xxxxxxxxxxxxxxxxxxxxxxxx") | Out-Null

$placeholderPara = $searchScope.Paragraphs(1).Range
$placeholderPara.InsertXML('<w:p><w:pPr><w:rPr><w:color w:val="0066CC"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0066CC"/></w:rPr><w:t>They were removed.</w:t></w:r></w:p>') | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new "_GoBack" bookmark in the middle of the word "unclear"
#    (…awkward or unc | lear sentences…) inside the final comment paragraph.
# ---------------------------------------------------------------------------
$splitPoint = $d.Content.Duplicate
$splitPoint.Find.Execute("Throughout the manuscript, there are several awkward or unc") | Out-Null
$insertAt = $d.Range($splitPoint.End, $splitPoint.End)
$d.Bookmarks.Add("_GoBack", $insertAt) | Out-Null
